$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The existing data row (row 22: YAHYA / 123 / 591EF2D4 / 0 / 103|100 / TRUE)
# is replaced by four new blank rows (22-25, styled like the rows above) and
# the data itself moves down to row 26 with updated values.

# Clear the old data out of row 22 so it becomes a blank row.
$ws.Range("A22:F22").ClearContents()

# Write the (moved + updated) data row at row 26.
# Leading "'" forces these to be stored as literal text (matching the
# original's shared-string cells) instead of being auto-typed as a number
# or boolean by Excel's input parser.
$ws.Range("A26").Value = "'YAHYA"
$ws.Range("B26").Value = "'YAHYA"
$ws.Range("C26").Value = "'EA4C7814"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = "'100"
$ws.Range("F26").Value = "'TRUE"

# Materialize rows 22-26 with the formatting (not the contents) of row 21,
# which already carries the plain un-bolded "data row" style used throughout
# rows 14-21. This both creates blank, properly-styled rows 22-25 and gives
# row 26 (the data row, just written above) the same style (must happen
# AFTER the values are written, otherwise the .Value assignment above resets
# the cell style back to default).
$ws.Range("A21:F21").Copy()
$ws.Range("A22:F26").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
